$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers in F1:H1, styled like the other headers (s="1")
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Boolean flag values for rows 2-19
$values = @(
    @($false, $false, $false),  # row 2
    @($false, $false, $false),  # row 3
    @($false, $false, $false),  # row 4
    @($true,  $false, $false),  # row 5
    @($true,  $true,  $true),   # row 6
    @($false, $false, $false),  # row 7
    @($false, $false, $false),  # row 8
    @($false, $false, $false),  # row 9
    @($false, $false, $false),  # row 10
    @($false, $false, $false),  # row 11
    @($false, $false, $false),  # row 12
    @($false, $false, $false),  # row 13
    @($false, $false, $false),  # row 14
    @($false, $false, $false),  # row 15
    @($false, $false, $false),  # row 16
    @($false, $false, $false),  # row 17
    @($false, $false, $false),  # row 18
    @($false, $false, $false)   # row 19
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i][0]
    $ws.Cells.Item($row, 7).Value = $values[$i][1]
    $ws.Cells.Item($row, 8).Value = $values[$i][2]
}
